$d = $word.ActiveDocument

function Get-ParaContaining($needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ------------------------------------------------------------------
# Change 1: paragraph "The table contains multiple categories..."
#   remove "film actor, " (first run text shrinks, rest of the
#   paragraph's runs are untouched content-wise)
# ------------------------------------------------------------------
$p1 = Get-ParaContaining("The table contains multiple categories")
$start1 = $p1.Range.Start
$full1 = $p1.Range.Text
$needle1 = "film actor, "
$idx1 = $full1.IndexOf($needle1)
$rng1 = $d.Range($start1 + $idx1, $start1 + $idx1 + $needle1.Length)
$rng1.Text = ""

# ------------------------------------------------------------------
# Change 2: paragraph "This tells me that..."
#   "Best Buy" -> "Blockbuster"
#   insert "automated like " right before "a Redbox"
# ------------------------------------------------------------------
$p2 = Get-ParaContaining("This tells me that")
$start2 = $p2.Range.Start
$full2 = $p2.Range.Text
$needle2 = "Best Buy"
$idx2 = $full2.IndexOf($needle2)
$rng2 = $d.Range($start2 + $idx2, $start2 + $idx2 + $needle2.Length)
$rng2.Text = "Blockbuster"

$p2b = Get-ParaContaining("This tells me that")
$start2b = $p2b.Range.Start
$full2b = $p2b.Range.Text
$needle3 = "a Redbox"
$idx3 = $full2b.IndexOf($needle3)
$rng3 = $d.Range($start2b + $idx3, $start2b + $idx3)
$rng3.Text = "automated like "

# ------------------------------------------------------------------
# Change 3: insert two new paragraphs (both indented ind=360) right
# after the "With your group: ... friend's business?" paragraph; the
# second one holds the new "We can add ID's..." sentence.
# ------------------------------------------------------------------
$p3 = Get-ParaContaining("friend")
$p3.Range.InsertParagraphAfter() | Out-Null

$p3 = Get-ParaContaining("friend")
$p3.Range.InsertParagraphAfter() | Out-Null

$p3 = Get-ParaContaining("friend")
$newPara1 = $p3.Next()
$newPara2 = $newPara1.Next()
$newPara2.Range.Text = "We can add ID" + [char]0x2019 + "s to every part as well as time data such as last update status. This makes it far easier to find data via the IDs and when events happened via the last update status."

# ------------------------------------------------------------------
# Mirror the re-paginated render the original commit captured: the
# trailing paragraph (just a manual line break) now starts a fresh
# page, so Word stamps a <w:lastRenderedPageBreak/> marker on it.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$insPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$insPoint.InsertXML("<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage' xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document><w:body><w:p><w:r><w:lastRenderedPageBreak/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")

Write-Output "Done"
